$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 9

# Row 4
$ws.Range("G4").Value = 3.8
$ws.Range("I4").Value = 2.1
$ws.Range("Q4").Value = 2.4
$ws.Range("R4").Value = 1.53
$ws.Range("X4").Value = 17
$ws.Range("AA4").Value = 34
$ws.Range("AB4").Value = 41
$ws.Range("AI4").Value = 9

# Row 8
$ws.Range("G8").Value = 2.55
$ws.Range("H8").Value = 3.6
$ws.Range("I8").Value = 2.4
$ws.Range("J8").Value = 3
$ws.Range("K8").Value = 2.27
$ws.Range("L8").Value = 2.9
$ws.Range("Q8").Value = 1.57
$ws.Range("R8").Value = 2.12
$ws.Range("U8").Value = 1.5
$ws.Range("V8").Value = 2.27
$ws.Range("W8").Value = 11.75
$ws.Range("Z8").Value = 29
$ws.Range("AA8").Value = 18.5
$ws.Range("AB8").Value = 22
$ws.Range("AE8").Value = 11.75
$ws.Range("AH8").Value = 11
$ws.Range("AI8").Value = 14
$ws.Range("AJ8").Value = 9.5
$ws.Range("AL8").Value = 17.5
$ws.Range("AM8").Value = 22
$ws.Range("AN8").Value = 4.7
$ws.Range("AO8").Value = 12.5
$ws.Range("AP8").Value = 17.5
$ws.Range("AQ8").Value = 50
$ws.Range("AU8").Value = 6.5
$ws.Range("AX8").Value = 12
$ws.Range("AY8").Value = 17.5
$ws.Range("BA8").Value = 65

# Row 14
$ws.Range("G14").Value = 2.47
$ws.Range("I14").Value = 2.65
$ws.Range("K14").Value = 2.07
$ws.Range("P14").Value = 2.9
$ws.Range("S14").Value = 1.42
$ws.Range("T14").Value = 2.65
$ws.Range("U14").Value = 1.82
$ws.Range("W14").Value = 7.7
$ws.Range("X14").Value = 12
$ws.Range("AE14").Value = 14.5
$ws.Range("AJ14").Value = 10
$ws.Range("AN14").Value = 4.4
$ws.Range("AT14").Value = 2.65
$ws.Range("AU14").Value = 7.1
$ws.Range("AV14").Value = 65
$ws.Range("AX14").Value = 14.5
$ws.Range("AY14").Value = 23
$ws.Range("BB14").Value = 300
